$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Team) -> short abbreviation codes, Column B (Year) -> numeric
# end-year of the season instead of the "YYYY-YY" string. Column C (Link)
# is unchanged.
$teams = @("GSW","CLE","LAL","GSW","CHI","BOS","LAL","CHI","CHI","SAS","BOS","LAL","DAL","MIA")
$years = @(2017,2016,2002,2015,1996,1986,1987,1991,1997,2014,2008,2009,2011,2013)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $teams[$i]
    $ws.Cells.Item($row, 2).Value = $years[$i]
}

$ws.Range("C23").Select()
